# Adding "Area" / "Atotal" computations to the discharge worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New headers (row 1) ----
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# ---- Row 2: Area for first segment, running Atotal, and the little J2:K2 summary ----
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# ---- D3:D8 turns into one shared "segment" formula group (previously plain, separate formulas) ----
$ws.Range("D3:D8").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# ---- G3 Area (first row of the Area column, not part of the G4:G15 shared group) ----
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# ---- G4:G15 shared "Area" formula group (extends two rows past the old data, to 13/14/15) ----
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Make the sheet's used range/selection match what was left selected in the source file.
[void]$ws.Range("J2:K2").Select()
